$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date (C1) ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45379

# --- "FPIEBP" sheet: re-prioritize "hard coal" (row 3) production/imports/exports
#     from (3,2,1) to (1,3,2), and move the active selection to E3 ---
$fpiebp = $wb.Worksheets.Item("FPIEBP")
$fpiebp.Activate()
$fpiebp.Range("B3").Value = 1
$fpiebp.Range("C3").Value = 3
$fpiebp.Range("D3").Value = 2
$fpiebp.Range("E3").Select()
